$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (col D) / 1h volume change (col E) snapshot values.
# Some "Price" values look like plain numbers (e.g. "312.28"); they are
# stored as text in the workbook, so they're written with a leading
# apostrophe to force text entry, then the cell style is reset to "Normal"
# so no stray quote-prefix style lingers on the cell.
$ws.Range("D2").Value = "42.138.89"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.312.85"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'312.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.08%  "
$ws.Range("E6").Value = "  +4.82%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").Value = "'0.0915"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "'8.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "'15.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.48%  "
$ws.Range("D16").Value = "2.658.20"
$ws.Range("D17").Value = "2.306.62"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "42.124.33"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "'7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").Value = "'74.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "'3.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.46%  "
$ws.Range("D23").Value = "'259.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").Value = "'9.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.17%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'10.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.69%  "
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "'35.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'162.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.77%  "
$ws.Range("E33").Value = "  -5.57%  "
$ws.Range("D34").Value = "'5.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.41%  "
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("E36").Value = "  +12.22%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").Value = "'2.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.41%  "
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").Value = "'98.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.53%  "
$ws.Range("E42").Value = "  -4.66%  "
$ws.Range("D43").Value = "'70.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "'0.229"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'12.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.44%  "
$ws.Range("D47").Value = "'111.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "'8.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").Value = "'74.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.73%  "
$ws.Range("E51").Value = "  -0.82%  "
